$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '89.599.40'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.043.19'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  -2.94%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.16'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '612.48'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -3.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.361'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -8.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.880'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +15.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.043.01'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -2.68%  '
$ws.Range("E11").Value = '  +21.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.187'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +4.85%  '
$ws.Range("E13").Value = '  -4.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.38'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.509.22'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.24'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -0.42%  '
$ws.Range("E17").Value = '  -3.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.045.62'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -3.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.31'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -2.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000218'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -5.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.38'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '424.01'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.23'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -1.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.02'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +2.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.36'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '83.82'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +1.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.62'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +0.74%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  +2.57%  '
$ws.Range("E30").Value = '  +0.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.25'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("B32").Value = 'dogwifhat'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.72'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -7.98%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '502.46'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -1.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.64'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -5.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.85'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +4.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.80'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -2.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.24'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -4.63%  '
$ws.Range("E38").Value = '  -10.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.24'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("E43").Value = '  +5.90%  '
$ws.Range("E44").Value = '  -1.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '146.77'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +1.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.24'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("E47").Value = '  +11.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.18'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +6.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '160.67'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -2.42%  '
$ws.Range("E50").Value = '  +2.29%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.71'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -2.44%  '
